# "Now with total maintenance time per day"
# Add Men / Menhours columns (C, D) and three new maintenance tasks
# (rows 5-7), plus a left/right "medium" grey border style used by the
# new Interval values in the lower block (B5:B7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- New header cells -------------------------------------------------
$ws.Range("C1").Value = "Men"
$ws.Range("D1").Value = "Menhours"

# ---- Existing task rows: Task Number becomes text, add Men/Menhours ---
$ws.Range("A2").Value = "0620000-00-01"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 0.2

$ws.Range("A3").Value = "0620000-00-02"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0.3

$ws.Range("A4").Value = "0620000-00-03"
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 0.3

# ---- New task rows ------------------------------------------------------
$ws.Range("A5").Value = "0620000-00-04"
$ws.Range("B5").Value = 10
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 1

$ws.Range("A6").Value = "0620000-00-05"
$ws.Range("B6").Value = 20
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 13

$ws.Range("A7").Value = "0620000-00-06"
$ws.Range("B7").Value = 21
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = 2

# Copy the existing Task Number / Interval cell formatting down onto the
# new rows (wraps text, right aligned, grey border) the same way it was
# applied to rows 2-4.
$ws.Range("A2:B4").Copy()
$ws.Range("A5:B7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# The Interval cells for the three new rows get a thinner look: only a
# left/right grey "medium" border (no top/bottom), matching the new
# maintenance-time block.
foreach ($r in 5..7) {
    $cell = $ws.Range("B$r")
    $left = $cell.Borders.Item(7)
    $left.Weight = -4138
    $left.Color = 13421772
    $right = $cell.Borders.Item(10)
    $right.Weight = -4138
    $right.Color = 13421772
}

# ---- Row heights --------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 14.4
$ws.Rows.Item(3).RowHeight = 15
$ws.Rows.Item(5).RowHeight = 14.4
$ws.Rows.Item(6).RowHeight = 15

# ---- Selection matches the saved workbook state --------------------------
$ws.Range("H21").Select()
